$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 1 header (unchanged visible text) ----
$ws.Range("A1").Value2 = "TabName"
$ws.Range("B1").Value2 = "query"
$ws.Range("C1").Value2 = "StatQuery"
$ws.Range("D1").Value2 = "dbExcel"
$ws.Range("E1").Value2 = "WebExcel"

# ---- Row 2 label ----
$ws.Range("A2").Value2 = "CasesTab"

# ---- Row 3 label (write first so shared-string order matches) ----
$ws.Range("A3").Value2 = "FilesTab"

# ---- Shared stat query (Trials / Cases / Files) used by both tabs ----
$statQuery = @"
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
   WHERE c.race = "AMERICAN_INDIAN_OR_ALASKA_NATIVE"
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files
"@
$ws.Range("C2").Value2 = $statQuery
$ws.Range("C2").WrapText = $true

# ---- Files-level filter query ----
$filesQuery = @"
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
OPTIONAL MATCH (f)-->(parent)
WITH f,a,ct,c,parent
  WHERE c.race = "AMERICAN_INDIAN_OR_ALASKA_NATIVE"
WITH
    f, parent, c, a, ct,
    ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
    toInteger(floor(log(f.file_size)/log(1024))) as i,
    2 as precision
WITH
    f, parent, c, a, ct,
    f.file_size /(1024^i) AS value,
    10^precision AS factor,
    units[i] as unit
WITH
    f, parent, c, a, ct, unit,
    round(factor * value)/factor AS size
RETURN DISTINCT
    f.file_name AS ``File Name``,
    head(labels(parent)) as Association,
    f.file_description AS Description,
    f.file_format AS ``File Format``,
    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    ct.clinical_trial_designation AS ``Trial Code``,
    a.arm_id AS Arm,
    c.case_id AS ``Case ID``
"@
$ws.Range("B3").Value2 = $filesQuery
$ws.Range("B3").WrapText = $true

# ---- Case-level filter query ----
$casesQuery = @"
MATCH (c:case)
 MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)
 MATCH (f:file)-[*]->(c)
WHERE c.race = "AMERICAN_INDIAN_OR_ALASKA_NATIVE"
RETURN DISTINCT
    c.case_id AS ``Case ID``,
     ct.clinical_trial_designation AS ``Trial Code``,
     a.arm_id AS Arm,
      a.arm_drug AS ``Arm Treatment``,
c.disease AS Diagnosis,
  c.gender AS Gender,
    c.race AS Race,
    c.ethnicity AS Ethnicity
"@
$ws.Range("B2").Value2 = $casesQuery
$ws.Range("B2").WrapText = $true

$ws.Range("C3").Value2 = $statQuery
$ws.Range("C3").WrapText = $true

# ---- Filenames (same on both rows) ----
$ws.Range("D2").Value2 = "TC01_Trials_Filter_Race-AmerIndAlask_Neo4jData.xlsx"
$ws.Range("E2").Value2 = "TC01_Trials_Filter_Race-AmerIndAlask_WebData.xlsx"
$ws.Range("D3").Value2 = "TC01_Trials_Filter_Race-AmerIndAlask_Neo4jData.xlsx"
$ws.Range("E3").Value2 = "TC01_Trials_Filter_Race-AmerIndAlask_WebData.xlsx"

# ---- Row heights ----
$ws.Rows.Item(2).RowHeight = 195
$ws.Rows.Item(3).RowHeight = 409.5

# ---- selection shown when the workbook is reopened ----
$ws.Range("A2:C3").Select() | Out-Null
